$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("http://www.emelevadores.com.br/", "contato@emelevadores.com.br;"),
    @("https://coteibem.sindiconet.com.br/fornecedores/manutencao-elevadores/sp/sao-paulo", "contato@coteibem.com.br;"),
    @("http://primac.com.br/", "comercial@primac.com.br;"),
    @("https://framartelelevadores.com.br/", "framartelelevadores@terra.com.br;"),
    @("https://continentalelevadores.com.br/", "continentalelevadores@protonmail.com;")
)

$startRow = 4
$endRow = $startRow + $data.Count - 1

$ws.Range("A3:B3").Copy() | Out-Null
$ws.Range("A$($startRow):B$($endRow)").PasteSpecial(-4122) | Out-Null

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
